$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 3).Value = 330296
$ws.Cells.Item(2, 4).Value = 420628070
$ws.Cells.Item(10, 3).Value = 119488
$ws.Cells.Item(10, 4).Value = 175068467
$ws.Cells.Item(11, 3).Value = 146
$ws.Cells.Item(11, 4).Value = 216013
$ws.Cells.Item(12, 3).Value = 61321
$ws.Cells.Item(12, 4).Value = 88493445
$ws.Cells.Item(16, 3).Value = 4052
$ws.Cells.Item(16, 4).Value = 5752637
$ws.Cells.Item(20, 3).Value = 7184
$ws.Cells.Item(20, 4).Value = 10030422
$ws.Cells.Item(22, 3).Value = 79411
$ws.Cells.Item(22, 4).Value = 98892614
$ws.Cells.Item(28, 3).Value = 33028
$ws.Cells.Item(28, 4).Value = 48338559
$ws.Cells.Item(30, 3).Value = 11789
$ws.Cells.Item(30, 4).Value = 16957239
$ws.Cells.Item(35, 3).Value = 1973
$ws.Cells.Item(35, 4).Value = 2785592
$ws.Cells.Item(36, 3).Value = 99434
$ws.Cells.Item(36, 4).Value = 124983216
$ws.Cells.Item(44, 3).Value = 45118
$ws.Cells.Item(44, 4).Value = 66111143
$ws.Cells.Item(46, 3).Value = 9406
$ws.Cells.Item(46, 4).Value = 13487123
$ws.Cells.Item(51, 3).Value = 2506
$ws.Cells.Item(51, 4).Value = 3503418
$ws.Cells.Item(52, 3).Value = 70758
$ws.Cells.Item(52, 4).Value = 88709482
$ws.Cells.Item(59, 3).Value = 28747
$ws.Cells.Item(59, 4).Value = 42158657
$ws.Cells.Item(62, 3).Value = 11437
$ws.Cells.Item(62, 4).Value = 16534544
$ws.Cells.Item(64, 3).Value = 1376
$ws.Cells.Item(64, 4).Value = 1923997
$ws.Cells.Item(68, 3).Value = 1591
$ws.Cells.Item(68, 4).Value = 2231363
$ws.Cells.Item(70, 3).Value = 20970
$ws.Cells.Item(70, 4).Value = 27458383
$ws.Cells.Item(74, 3).Value = 7749
$ws.Cells.Item(74, 4).Value = 11349148
$ws.Cells.Item(76, 3).Value = 5244
$ws.Cells.Item(76, 4).Value = 7615294
$ws.Cells.Item(79, 3).Value = 144456
$ws.Cells.Item(79, 4).Value = 179976756
$ws.Cells.Item(85, 3).Value = 64793
$ws.Cells.Item(85, 4).Value = 94954309
$ws.Cells.Item(88, 3).Value = 30492
$ws.Cells.Item(88, 4).Value = 44108470
$ws.Cells.Item(91, 3).Value = 3026
$ws.Cells.Item(91, 4).Value = 4276010
$ws.Cells.Item(92, 3).Value = 34744
$ws.Cells.Item(92, 4).Value = 47110522
$ws.Cells.Item(96, 3).Value = 8419
$ws.Cells.Item(96, 4).Value = 12376097
$ws.Cells.Item(98, 3).Value = 7806
$ws.Cells.Item(98, 4).Value = 11330807
$ws.Cells.Item(100, 3).Value = 555
$ws.Cells.Item(100, 4).Value = 787656
$ws.Cells.Item(101, 3).Value = 525
$ws.Cells.Item(101, 4).Value = 758050
$ws.Cells.Item(102, 3).Value = 11789
$ws.Cells.Item(102, 4).Value = 19125984
$ws.Cells.Item(104, 3).Value = 2857
$ws.Cells.Item(104, 4).Value = 4980170
$ws.Cells.Item(106, 3).Value = 3937
$ws.Cells.Item(106, 4).Value = 6892269
$ws.Cells.Item(109, 3).Value = 227
$ws.Cells.Item(109, 4).Value = 370530
$ws.Cells.Item(110, 3).Value = 145523
$ws.Cells.Item(110, 4).Value = 179972407
$ws.Cells.Item(116, 3).Value = 53812
$ws.Cells.Item(116, 4).Value = 78857461
$ws.Cells.Item(118, 3).Value = 28073
$ws.Cells.Item(118, 4).Value = 40671950
$ws.Cells.Item(122, 3).Value = 2414
$ws.Cells.Item(122, 4).Value = 3395994
$ws.Cells.Item(124, 3).Value = 543135
$ws.Cells.Item(124, 4).Value = 717740829
$ws.Cells.Item(129, 3).Value = 1415
$ws.Cells.Item(129, 4).Value = 2097214
$ws.Cells.Item(131, 3).Value = 215477
$ws.Cells.Item(131, 4).Value = 316715341
$ws.Cells.Item(134, 3).Value = 193946
$ws.Cells.Item(134, 4).Value = 282040188
$ws.Cells.Item(137, 3).Value = 2891
$ws.Cells.Item(137, 4).Value = 4059072
$ws.Cells.Item(140, 3).Value = 6874
$ws.Cells.Item(140, 4).Value = 9698876
$ws.Cells.Item(143, 3).Value = 46259
$ws.Cells.Item(143, 4).Value = 61739950
$ws.Cells.Item(147, 3).Value = 6
$ws.Cells.Item(147, 4).Value = 9000
$ws.Cells.Item(149, 3).Value = 14444
$ws.Cells.Item(149, 4).Value = 21170094
$ws.Cells.Item(150, 3).Value = 3880
$ws.Cells.Item(150, 4).Value = 5594931
$ws.Cells.Item(153, 3).Value = 407
$ws.Cells.Item(153, 4).Value = 585716
$ws.Cells.Item(155, 3).Value = 417
$ws.Cells.Item(155, 4).Value = 588313
$ws.Cells.Item(156, 3).Value = 18250
$ws.Cells.Item(156, 4).Value = 24130080
$ws.Cells.Item(160, 3).Value = 7452
$ws.Cells.Item(160, 4).Value = 10846728
$ws.Cells.Item(162, 3).Value = 5214
$ws.Cells.Item(162, 4).Value = 7505821
$ws.Cells.Item(165, 3).Value = 281
$ws.Cells.Item(165, 4).Value = 401664
$ws.Cells.Item(167, 3).Value = 21587
$ws.Cells.Item(167, 4).Value = 38346680
$ws.Cells.Item(168, 3).Value = 2295
$ws.Cells.Item(168, 4).Value = 4067571
$ws.Cells.Item(169, 3).Value = 301
$ws.Cells.Item(169, 4).Value = 522089
$ws.Cells.Item(173, 3).Value = 89977
$ws.Cells.Item(173, 4).Value = 112364889
$ws.Cells.Item(177, 3).Value = 15
$ws.Cells.Item(177, 4).Value = 21889
$ws.Cells.Item(180, 3).Value = 34532
$ws.Cells.Item(180, 4).Value = 50632229
$ws.Cells.Item(181, 3).Value = 29
$ws.Cells.Item(181, 4).Value = 43500
$ws.Cells.Item(182, 3).Value = 13410
$ws.Cells.Item(182, 4).Value = 19375414
$ws.Cells.Item(184, 3).Value = 1269
$ws.Cells.Item(184, 4).Value = 1776027
$ws.Cells.Item(186, 3).Value = 1755
$ws.Cells.Item(186, 4).Value = 2463929
$ws.Cells.Item(188, 3).Value = 244403
$ws.Cells.Item(188, 4).Value = 303594093
$ws.Cells.Item(190, 3).Value = 174
$ws.Cells.Item(190, 4).Value = 251236
$ws.Cells.Item(196, 3).Value = 88215
$ws.Cells.Item(196, 4).Value = 129289770
$ws.Cells.Item(199, 3).Value = 33911
$ws.Cells.Item(199, 4).Value = 48817678
$ws.Cells.Item(202, 3).Value = 5194
$ws.Cells.Item(202, 4).Value = 7395005
$ws.Cells.Item(205, 3).Value = 5212
$ws.Cells.Item(205, 4).Value = 7217775
$ws.Cells.Item(208, 3).Value = 270954
$ws.Cells.Item(208, 4).Value = 335204816
$ws.Cells.Item(210, 3).Value = 259
$ws.Cells.Item(210, 4).Value = 370039
$ws.Cells.Item(215, 3).Value = 628
$ws.Cells.Item(215, 4).Value = 914878
$ws.Cells.Item(217, 3).Value = 97125
$ws.Cells.Item(217, 4).Value = 142083087
$ws.Cells.Item(218, 3).Value = 97
$ws.Cells.Item(218, 4).Value = 144699
$ws.Cells.Item(220, 3).Value = 53039
$ws.Cells.Item(220, 4).Value = 76653644
$ws.Cells.Item(226, 3).Value = 6216
$ws.Cells.Item(226, 4).Value = 8612618
$ws.Cells.Item(229, 3).Value = 109243
$ws.Cells.Item(229, 4).Value = 136530613
$ws.Cells.Item(231, 3).Value = 78
$ws.Cells.Item(231, 4).Value = 112013
$ws.Cells.Item(236, 3).Value = 50469
$ws.Cells.Item(236, 4).Value = 73928892
$ws.Cells.Item(237, 3).Value = 40
$ws.Cells.Item(237, 4).Value = 57711
$ws.Cells.Item(238, 3).Value = 12873
$ws.Cells.Item(238, 4).Value = 18517067
$ws.Cells.Item(242, 3).Value = 2679
$ws.Cells.Item(242, 4).Value = 3753506
$ws.Cells.Item(243, 3).Value = 265590
$ws.Cells.Item(243, 4).Value = 335314089
$ws.Cells.Item(251, 3).Value = 98002
$ws.Cells.Item(251, 4).Value = 143590654
$ws.Cells.Item(253, 3).Value = 5
$ws.Cells.Item(253, 4).Value = 7500
$ws.Cells.Item(254, 3).Value = 67385
$ws.Cells.Item(254, 4).Value = 97677937
$ws.Cells.Item(256, 3).Value = 2453
$ws.Cells.Item(256, 4).Value = 3460224
$ws.Cells.Item(259, 3).Value = 4959
$ws.Cells.Item(259, 4).Value = 6963977
